$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.066.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.910.62'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.22%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.88%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.14'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4819'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.75%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3817'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07358'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9352'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.84'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07783'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.895.10'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.509'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.626'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.79'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008832'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '28.100.29'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.03%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.177'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.119.58'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.79%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.09'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.921'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.56'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.110'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +4.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.63'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.960'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08959'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.300'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.24%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7776'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.682'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.657'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.30%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05312'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5479'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.985'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.021'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.508'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.66'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4845'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '108.49'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.005'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.655'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.11'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06086'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.04%  '
